$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 -> "Enemy 2"
$ws.Range("A2").Value = "Enemy 2"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 14
$ws.Range("G2").Value = 25

# Update row 3 -> "Enemy 1"
$ws.Range("A3").Value = "Enemy 1"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 14
$ws.Range("G3").Value = 25

# Delete rows 4-7
$ws.Range("A4:K7").EntireRow.Delete()
